$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.526.53"
$ws.Range("D3").Value = "1.830.18"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.88"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3609"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07220"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8609"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.822.84"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.393"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.494"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06923"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.55"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008906"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.38"
$ws.Range("D21").Value = "27.626.07"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.149"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("E23").Value = "  +4.43%  "
$ws.Range("D24").Value = "2.083.13"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.985"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.26"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.151"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.60"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.796"
$ws.Range("E30").Value = "  -5.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08934"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7479"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.552"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.965"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.084"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05260"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01927"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.800"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5080"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1655"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.359"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.373"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.83"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.43"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06467"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4683"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9996"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.618"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.92"
$ws.Range("E51").Value = "  -0.98%  "
